# Auto-generated edit script applying numeric value updates to the
# Leve profit-tracking tables (columns H-N) across all 8 job sheets,
# as produced by the scheduled price-refresh runner.
$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(6, 8).Value = 0  # H6: 63 -> 0
$ws.Cells.Item(6, 9).Value = 0  # I6: 63 -> 0
$ws.Cells.Item(6, 11).Value = 0  # K6: 189 -> 0
$ws.Cells.Item(6, 13).ClearContents()  # M6: -77 -> (removed)
$ws.Cells.Item(9, 8).Value = 1523.8889  # H9: 1172.9166 -> 1523.8889
$ws.Cells.Item(9, 9).Value = 1984.3334  # I9: 1362.8889 -> 1984.3334
$ws.Cells.Item(9, 11).Value = 1984.3334  # K9: 1362.8889 -> 1984.3334
$ws.Cells.Item(9, 13).Value = -1815.3334  # M9: -1193.8889 -> -1815.3334
$ws.Cells.Item(18, 8).Value = 3999.5  # H18: 1692.625 -> 3999.5
$ws.Cells.Item(18, 9).Value = 5000  # I18: 1506 -> 5000
$ws.Cells.Item(18, 11).Value = 5000  # K18: 1506 -> 5000
$ws.Cells.Item(18, 13).Value = -4716  # M18: -1222 -> -4716
$ws.Cells.Item(32, 8).Value = 4458  # H32: 4022.625 -> 4458
$ws.Cells.Item(32, 9).Value = 1445.5  # I32: 1260.3334 -> 1445.5
$ws.Cells.Item(32, 10).Value = 6466.3335  # J32: 5680 -> 6466.3335
$ws.Cells.Item(32, 11).Value = 1445.5  # K32: 1260.3334 -> 1445.5
$ws.Cells.Item(32, 12).Value = 6466.3335  # L32: 5680 -> 6466.3335
$ws.Cells.Item(32, 13).Value = -1119.5  # M32: -934.3334 -> -1119.5
$ws.Cells.Item(32, 14).Value = -7118.3335  # N32: -6332 -> -7118.3335
$ws.Cells.Item(40, 8).Value = 1653  # H40: 1671.0714 -> 1653
$ws.Cells.Item(40, 10).Value = 1779.8  # J40: 1874.75 -> 1779.8
$ws.Cells.Item(40, 12).Value = 1779.8  # L40: 1874.75 -> 1779.8
$ws.Cells.Item(40, 14).Value = -2129.8  # N40: -2224.75 -> -2129.8
$ws.Cells.Item(42, 8).Value = 3645.889  # H42: 4644.8335 -> 3645.889
$ws.Cells.Item(42, 9).Value = 333.66666  # I42: 310.6 -> 333.66666
$ws.Cells.Item(42, 10).Value = 10270.333  # J42: 26316 -> 10270.333
$ws.Cells.Item(42, 11).Value = 1000.99998  # K42: 931.8000000000001 -> 1000.99998
$ws.Cells.Item(42, 12).Value = 30810.999  # L42: 78948 -> 30810.999
$ws.Cells.Item(42, 13).Value = -770.9999799999999  # M42: -701.8000000000001 -> -770.9999799999999
$ws.Cells.Item(42, 14).Value = -31270.999  # N42: -79408 -> -31270.999
$ws.Cells.Item(51, 8).Value = 4870.727  # H51: 4707.154 -> 4870.727
$ws.Cells.Item(51, 9).Value = 3797.5715  # I51: 3799.75 -> 3797.5715
$ws.Cells.Item(51, 10).Value = 6748.75  # J51: 6159 -> 6748.75
$ws.Cells.Item(51, 11).Value = 3797.5715  # K51: 3799.75 -> 3797.5715
$ws.Cells.Item(51, 12).Value = 6748.75  # L51: 6159 -> 6748.75
$ws.Cells.Item(51, 13).Value = -3313.5715  # M51: -3315.75 -> -3313.5715
$ws.Cells.Item(51, 14).Value = -7716.75  # N51: -7127 -> -7716.75
$ws.Cells.Item(64, 8).Value = 5500  # H64: 5498.5 -> 5500
$ws.Cells.Item(64, 10).Value = 5500  # J64: 5498.5 -> 5500
$ws.Cells.Item(64, 12).Value = 5500  # L64: 5498.5 -> 5500
$ws.Cells.Item(64, 14).Value = -5996  # N64: -5994.5 -> -5996
$ws.Cells.Item(67, 8).Value = 5500  # H67: 5498.5 -> 5500
$ws.Cells.Item(67, 10).Value = 5500  # J67: 5498.5 -> 5500
$ws.Cells.Item(67, 12).Value = 5500  # L67: 5498.5 -> 5500
$ws.Cells.Item(67, 14).Value = -7216  # N67: -7214.5 -> -7216
$ws.Cells.Item(69, 8).Value = 21000.646  # H69: 21294.766 -> 21000.646
$ws.Cells.Item(69, 10).Value = 21250.875  # J69: 21563.375 -> 21250.875
$ws.Cells.Item(69, 12).Value = 63752.625  # L69: 64690.125 -> 63752.625
$ws.Cells.Item(69, 14).Value = -65500.625  # N69: -66438.125 -> -65500.625
$ws.Cells.Item(72, 8).Value = 21000.646  # H72: 21294.766 -> 21000.646
$ws.Cells.Item(72, 10).Value = 21250.875  # J72: 21563.375 -> 21250.875
$ws.Cells.Item(72, 12).Value = 191257.875  # L72: 194070.375 -> 191257.875
$ws.Cells.Item(72, 14).Value = -199993.875  # N72: -202806.375 -> -199993.875
$ws.Cells.Item(80, 8).Value = 4307.9165  # H80: 4627.091 -> 4307.9165
$ws.Cells.Item(80, 10).Value = 4041.8572  # J80: 4582.6665 -> 4041.8572
$ws.Cells.Item(80, 12).Value = 12125.5716  # L80: 13747.9995 -> 12125.5716
$ws.Cells.Item(80, 14).Value = -14121.5716  # N80: -15743.9995 -> -14121.5716
$ws.Cells.Item(83, 8).Value = 4307.9165  # H83: 4627.091 -> 4307.9165
$ws.Cells.Item(83, 10).Value = 4041.8572  # J83: 4582.6665 -> 4041.8572
$ws.Cells.Item(83, 12).Value = 36376.7148  # L83: 41243.9985 -> 36376.7148
$ws.Cells.Item(83, 14).Value = -46360.7148  # N83: -51227.9985 -> -46360.7148
$ws.Cells.Item(88, 8).Value = 1347  # H88: 0 -> 1347
$ws.Cells.Item(88, 9).Value = 1295  # I88: 0 -> 1295
$ws.Cells.Item(88, 10).Value = 1399  # J88: 0 -> 1399
$ws.Cells.Item(88, 11).Value = 1295  # K88: 0 -> 1295
$ws.Cells.Item(88, 12).Value = 1399  # L88: 0 -> 1399
$ws.Cells.Item(88, 13).Value = -889  # M88: None -> -889
$ws.Cells.Item(88, 14).Value = -2211  # N88: None -> -2211
$ws.Cells.Item(91, 8).Value = 1347  # H91: 0 -> 1347
$ws.Cells.Item(91, 9).Value = 1295  # I91: 0 -> 1295
$ws.Cells.Item(91, 10).Value = 1399  # J91: 0 -> 1399
$ws.Cells.Item(91, 11).Value = 1295  # K91: 0 -> 1295
$ws.Cells.Item(91, 12).Value = 1399  # L91: 0 -> 1399
$ws.Cells.Item(91, 13).Value = 109  # M91: None -> 109
$ws.Cells.Item(91, 14).Value = -4207  # N91: None -> -4207
$ws.Cells.Item(113, 8).Value = 2099  # H113: 1740.8 -> 2099
$ws.Cells.Item(113, 9).Value = 798  # I113: 800 -> 798
$ws.Cells.Item(113, 10).Value = 2749.5  # J113: 1976 -> 2749.5
$ws.Cells.Item(113, 11).Value = 798  # K113: 800 -> 798
$ws.Cells.Item(113, 12).Value = 2749.5  # L113: 1976 -> 2749.5
$ws.Cells.Item(113, 13).Value = 2456  # M113: 2454 -> 2456
$ws.Cells.Item(113, 14).Value = -9257.5  # N113: -8484 -> -9257.5
$ws.Cells.Item(116, 8).Value = 0  # H116: 1899 -> 0
$ws.Cells.Item(116, 10).Value = 0  # J116: 1899 -> 0
$ws.Cells.Item(116, 12).Value = 0  # L116: 1899 -> 0
$ws.Cells.Item(116, 14).ClearContents()  # N116: -8783 -> (removed)
$ws.Cells.Item(129, 8).Value = 2694.0952  # H129: 2671.238 -> 2694.0952
$ws.Cells.Item(129, 9).Value = 1098.6666  # I129: 738.6667 -> 1098.6666
$ws.Cells.Item(129, 10).Value = 2960  # J129: 2993.3333 -> 2960
$ws.Cells.Item(129, 11).Value = 3295.9998  # K129: 2216.0001 -> 3295.9998
$ws.Cells.Item(129, 12).Value = 8880  # L129: 8979.999899999999 -> 8880
$ws.Cells.Item(129, 13).Value = 1704.0002  # M129: 2783.9999 -> 1704.0002
$ws.Cells.Item(129, 14).Value = -18880  # N129: -18979.9999 -> -18880
$ws.Cells.Item(141, 8).Value = 3421  # H141: 3243.375 -> 3421
$ws.Cells.Item(141, 9).Value = 2624.25  # I141: 2499.4 -> 2624.25
$ws.Cells.Item(141, 11).Value = 7872.75  # K141: 7498.200000000001 -> 7872.75
$ws.Cells.Item(141, 13).Value = -2692.75  # M141: -2318.200000000001 -> -2692.75

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(94, 8).Value = 22330  # H94: 34886.668 -> 22330
$ws.Cells.Item(94, 10).Value = 22330  # J94: 34886.668 -> 22330
$ws.Cells.Item(94, 12).Value = 22330  # L94: 34886.668 -> 22330
$ws.Cells.Item(94, 14).Value = -24132  # N94: -36688.668 -> -24132
$ws.Cells.Item(97, 8).Value = 1900  # H97: 0 -> 1900
$ws.Cells.Item(97, 10).Value = 1900  # J97: 0 -> 1900
$ws.Cells.Item(97, 12).Value = 1900  # L97: 0 -> 1900
$ws.Cells.Item(97, 14).Value = -2892  # N97: None -> -2892

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(86, 8).Value = 1829.2222  # H86: 2083.1667 -> 1829.2222
$ws.Cells.Item(86, 9).Value = 1745.375  # I86: 1999.8 -> 1745.375
$ws.Cells.Item(86, 11).Value = 1745.375  # K86: 1999.8 -> 1745.375
$ws.Cells.Item(86, 13).Value = -622.375  # M86: -876.8 -> -622.375
$ws.Cells.Item(89, 8).Value = 1829.2222  # H89: 2083.1667 -> 1829.2222
$ws.Cells.Item(89, 9).Value = 1745.375  # I89: 1999.8 -> 1745.375
$ws.Cells.Item(89, 11).Value = 8726.875  # K89: 9999 -> 8726.875
$ws.Cells.Item(89, 13).Value = -3110.875  # M89: -4383 -> -3110.875
$ws.Cells.Item(135, 8).Value = 44875.883  # H135: 45000 -> 44875.883
$ws.Cells.Item(135, 10).Value = 44875.883  # J135: 45000 -> 44875.883
$ws.Cells.Item(135, 12).Value = 44875.883  # L135: 45000 -> 44875.883
$ws.Cells.Item(135, 14).Value = -55015.883  # N135: -55140 -> -55015.883

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(58, 8).Value = 1713.0667  # H58: 1849.5714 -> 1713.0667
$ws.Cells.Item(58, 9).Value = 1561.3077  # I58: 1509.6 -> 1561.3077
$ws.Cells.Item(58, 11).Value = 1561.3077  # K58: 1509.6 -> 1561.3077
$ws.Cells.Item(58, 13).Value = -1358.3077  # M58: -1306.6 -> -1358.3077
$ws.Cells.Item(99, 8).Value = 0  # H99: 2449 -> 0
$ws.Cells.Item(99, 10).Value = 0  # J99: 2449 -> 0
$ws.Cells.Item(99, 12).Value = 0  # L99: 2449 -> 0
$ws.Cells.Item(99, 14).ClearContents()  # N99: -5445 -> (removed)
$ws.Cells.Item(126, 8).Value = 0  # H126: 2449 -> 0
$ws.Cells.Item(126, 10).Value = 0  # J126: 2449 -> 0
$ws.Cells.Item(126, 12).Value = 0  # L126: 7347 -> 0
$ws.Cells.Item(126, 14).ClearContents()  # N126: -12287 -> (removed)
$ws.Cells.Item(132, 8).Value = 4240.857  # H132: 4631.1665 -> 4240.857
$ws.Cells.Item(132, 9).Value = 4997.8  # I132: 5772.5 -> 4997.8
$ws.Cells.Item(132, 11).Value = 14993.4  # K132: 17317.5 -> 14993.4
$ws.Cells.Item(132, 13).Value = -12463.4  # M132: -14787.5 -> -12463.4
$ws.Cells.Item(134, 8).Value = 2040.6666  # H134: 2120 -> 2040.6666
$ws.Cells.Item(134, 9).Value = 2040.6666  # I134: 2120 -> 2040.6666
$ws.Cells.Item(134, 11).Value = 6121.9998  # K134: 6360 -> 6121.9998
$ws.Cells.Item(134, 13).Value = -3586.9998  # M134: -3825 -> -3586.9998
$ws.Cells.Item(136, 8).Value = 1713.0667  # H136: 1849.5714 -> 1713.0667
$ws.Cells.Item(136, 9).Value = 1561.3077  # I136: 1509.6 -> 1561.3077
$ws.Cells.Item(136, 11).Value = 4683.9231  # K136: 4528.799999999999 -> 4683.9231
$ws.Cells.Item(136, 13).Value = -2133.9231  # M136: -1978.799999999999 -> -2133.9231

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(68, 8).Value = 1079  # H68: 1033 -> 1079
$ws.Cells.Item(68, 9).Value = 346.5  # I68: 350 -> 346.5
$ws.Cells.Item(68, 10).Value = 1262.125  # J68: 1228.1428 -> 1262.125
$ws.Cells.Item(68, 11).Value = 1039.5  # K68: 1050 -> 1039.5
$ws.Cells.Item(68, 12).Value = 3786.375  # L68: 3684.4284 -> 3786.375
$ws.Cells.Item(68, 13).Value = -228.5  # M68: -239 -> -228.5
$ws.Cells.Item(68, 14).Value = -5408.375  # N68: -5306.428400000001 -> -5408.375
$ws.Cells.Item(71, 8).Value = 1079  # H71: 1033 -> 1079
$ws.Cells.Item(71, 9).Value = 346.5  # I71: 350 -> 346.5
$ws.Cells.Item(71, 10).Value = 1262.125  # J71: 1228.1428 -> 1262.125
$ws.Cells.Item(71, 11).Value = 3118.5  # K71: 3150 -> 3118.5
$ws.Cells.Item(71, 12).Value = 11359.125  # L71: 11053.2852 -> 11359.125
$ws.Cells.Item(71, 13).Value = 937.5  # M71: 906 -> 937.5
$ws.Cells.Item(71, 14).Value = -19471.125  # N71: -19165.2852 -> -19471.125
$ws.Cells.Item(92, 8).Value = 683.6  # H92: 662.9524 -> 683.6
$ws.Cells.Item(92, 9).Value = 549.5  # I92: 489.6 -> 549.5
$ws.Cells.Item(92, 11).Value = 1648.5  # K92: 1468.8 -> 1648.5
$ws.Cells.Item(92, 13).Value = -400.5  # M92: -220.8000000000002 -> -400.5
$ws.Cells.Item(98, 8).Value = 1363.5555  # H98: 1452.6666 -> 1363.5555
$ws.Cells.Item(98, 9).Value = 999  # I98: 943.5 -> 999
$ws.Cells.Item(98, 10).Value = 1409.125  # J98: 1598.1428 -> 1409.125
$ws.Cells.Item(98, 11).Value = 2997  # K98: 2830.5 -> 2997
$ws.Cells.Item(98, 12).Value = 4227.375  # L98: 4794.428400000001 -> 4227.375
$ws.Cells.Item(98, 13).Value = -1499  # M98: -1332.5 -> -1499
$ws.Cells.Item(98, 14).Value = -7223.375  # N98: -7790.428400000001 -> -7223.375
$ws.Cells.Item(113, 8).Value = 1692.1428  # H113: 2114.5 -> 1692.1428
$ws.Cells.Item(113, 10).Value = 1600  # J113: 2162.5 -> 1600
$ws.Cells.Item(113, 12).Value = 4800  # L113: 6487.5 -> 4800
$ws.Cells.Item(113, 14).Value = -9140  # N113: -10827.5 -> -9140

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(80, 8).Value = 3871.7273  # H80: 3949.3 -> 3871.7273
$ws.Cells.Item(80, 9).Value = 3849.25  # I80: 3849.75 -> 3849.25
$ws.Cells.Item(80, 10).Value = 3884.5715  # J80: 4015.6667 -> 3884.5715
$ws.Cells.Item(80, 11).Value = 3849.25  # K80: 3849.75 -> 3849.25
$ws.Cells.Item(80, 12).Value = 3884.5715  # L80: 4015.6667 -> 3884.5715
$ws.Cells.Item(80, 13).Value = -2851.25  # M80: -2851.75 -> -2851.25
$ws.Cells.Item(80, 14).Value = -5880.5715  # N80: -6011.6667 -> -5880.5715
$ws.Cells.Item(83, 8).Value = 3871.7273  # H83: 3949.3 -> 3871.7273
$ws.Cells.Item(83, 9).Value = 3849.25  # I83: 3849.75 -> 3849.25
$ws.Cells.Item(83, 10).Value = 3884.5715  # J83: 4015.6667 -> 3884.5715
$ws.Cells.Item(83, 11).Value = 19246.25  # K83: 19248.75 -> 19246.25
$ws.Cells.Item(83, 12).Value = 19422.8575  # L83: 20078.3335 -> 19422.8575
$ws.Cells.Item(83, 13).Value = -14254.25  # M83: -14256.75 -> -14254.25
$ws.Cells.Item(83, 14).Value = -29406.8575  # N83: -30062.3335 -> -29406.8575
$ws.Cells.Item(102, 8).Value = 3696  # H102: 3519.6 -> 3696
$ws.Cells.Item(102, 10).Value = 0  # J102: 2814 -> 0
$ws.Cells.Item(102, 12).Value = 0  # L102: 2814 -> 0
$ws.Cells.Item(102, 14).ClearContents()  # N102: -6058 -> (removed)
$ws.Cells.Item(113, 8).Value = 799  # H113: 870.4286 -> 799
$ws.Cells.Item(113, 9).Value = 799  # I113: 870.4286 -> 799
$ws.Cells.Item(113, 11).Value = 799  # K113: 870.4286 -> 799
$ws.Cells.Item(113, 13).Value = 1371  # M113: 1299.5714 -> 1371

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(16, 8).Value = 788.4  # H16: 573.9231 -> 788.4
$ws.Cells.Item(16, 9).Value = 788.4  # I16: 588.5 -> 788.4
$ws.Cells.Item(16, 10).Value = 0  # J16: 399 -> 0
$ws.Cells.Item(16, 11).Value = 788.4  # K16: 588.5 -> 788.4
$ws.Cells.Item(16, 12).Value = 0  # L16: 399 -> 0
$ws.Cells.Item(16, 13).Value = -618.4  # M16: -418.5 -> -618.4
$ws.Cells.Item(16, 14).ClearContents()  # N16: -739 -> (removed)
$ws.Cells.Item(40, 8).Value = 2528  # H40: 2770.2856 -> 2528
$ws.Cells.Item(40, 9).Value = 2528  # I40: 2770.2856 -> 2528
$ws.Cells.Item(40, 11).Value = 2528  # K40: 2770.2856 -> 2528
$ws.Cells.Item(40, 13).Value = -2392  # M40: -2634.2856 -> -2392
$ws.Cells.Item(46, 8).Value = 1576.8  # H46: 1646.25 -> 1576.8
$ws.Cells.Item(46, 9).Value = 1696.25  # I46: 1828.6666 -> 1696.25
$ws.Cells.Item(46, 11).Value = 1696.25  # K46: 1828.6666 -> 1696.25
$ws.Cells.Item(46, 13).Value = -1508.25  # M46: -1640.6666 -> -1508.25
$ws.Cells.Item(122, 8).Value = 5749.2085  # H122: 5908.591 -> 5749.2085
$ws.Cells.Item(122, 9).Value = 4749.1875  # I122: 4856.7856 -> 4749.1875
$ws.Cells.Item(122, 11).Value = 14247.5625  # K122: 14570.3568 -> 14247.5625
$ws.Cells.Item(122, 13).Value = -11797.5625  # M122: -12120.3568 -> -11797.5625

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(46, 8).Value = 38857.2  # H46: 45412.223 -> 38857.2
$ws.Cells.Item(46, 10).Value = 38857.2  # J46: 45412.223 -> 38857.2
$ws.Cells.Item(46, 12).Value = 38857.2  # L46: 45412.223 -> 38857.2
$ws.Cells.Item(46, 14).Value = -39319.2  # N46: -45874.223 -> -39319.2
$ws.Cells.Item(100, 8).Value = 16668260  # H100: 12501861 -> 16668260
$ws.Cells.Item(100, 9).Value = 20001468  # I100: 16668407 -> 20001468
$ws.Cells.Item(100, 10).Value = 2222  # J100: 2221.5 -> 2222
$ws.Cells.Item(100, 11).Value = 40002936  # K100: 33336814 -> 40002936
$ws.Cells.Item(100, 12).Value = 4444  # L100: 4443 -> 4444
$ws.Cells.Item(100, 13).Value = -40002395  # M100: -33336273 -> -40002395
$ws.Cells.Item(100, 14).Value = -5526  # N100: -5525 -> -5526
$ws.Cells.Item(134, 8).Value = 38857.2  # H134: 45412.223 -> 38857.2
$ws.Cells.Item(134, 10).Value = 38857.2  # J134: 45412.223 -> 38857.2
$ws.Cells.Item(134, 12).Value = 116571.6  # L134: 136236.669 -> 116571.6
$ws.Cells.Item(134, 14).Value = -121641.6  # N134: -141306.669 -> -121641.6
